# Update the "Sites_networks" sheet: split the old "Insight_internal" /
# 172.16.0.0/12 row into the new, narrower supernets (172.16.0.0/13 plus
# the individual /16s that are no longer covered by 172.16.0.0/13), and
# drop the now-redundant "Insight_others" / summarized ranges.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sites_networks")

$ws.Range("A2").Value = "Insight_lab"
$ws.Range("B2").Value = "172.28.0.0/16"

$ws.Range("A3").Value  = "Insight_internal"
$ws.Range("B3").Value  = "172.16.0.0/13"
$ws.Range("A4").Value  = "Insight_internal"
$ws.Range("B4").Value  = "172.24.0.0/16"
$ws.Range("A5").Value  = "Insight_internal"
$ws.Range("B5").Value  = "172.25.0.0/16"
$ws.Range("A6").Value  = "Insight_internal"
$ws.Range("B6").Value  = "172.26.0.0/16"
$ws.Range("A7").Value  = "Insight_internal"
$ws.Range("B7").Value  = "172.27.0.0/16"
$ws.Range("A8").Value  = "Insight_internal"
$ws.Range("B8").Value  = "172.29.0.0/16"
$ws.Range("A9").Value  = "Insight_internal"
$ws.Range("B9").Value  = "172.30.0.0/16"
$ws.Range("A10").Value = "Insight_internal"
$ws.Range("B10").Value = "172.31.0.0/16"

$ws.Range("A11").Value = "Insight_Azure"
$ws.Range("B11").Value = "10.0.1.0/24"

# Add a warning to the "networks" column comment about the lack of
# overlapping-summarization support.
$comment = $ws.Range("B1").Comment
$existingText = $comment.Text()
[void]$comment.Text($existingText + "`nIMPORTANT: DOES NOT support overlaping summarization you MUST exclude manually the range")

# Make "Sites_networks" the active sheet/tab, with E11 selected (the last
# data row just entered).
[void]$ws.Activate()
[void]$ws.Range("E11").Select()
